# Commit: "Removed unused analyses and removed side-scoring"
#
# 1. hypothesisTests sheet: refresh the computed stats (t, p, df moved to a
#    fixed 173 and the B:D/E number formats got simplified), then move the
#    selection.
# 2. Add a new "plex" worksheet (same layout) with the per-motivator
#    hypothesis tests, then hand focus back to hypothesisTests.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("hypothesisTests")

# ---------------------------------------------------------------------------
# 1. Update the existing hypothesisTests numbers (rows 2-7, cols B:F)
# ---------------------------------------------------------------------------

$hypoRows = @(
    @(2, "4.8470588235294096",    "5.0999999999999996",   "-1.1297195597313301",   "0.26015835458148001",   "173"),
    @(3, "5.1529411764705797",    "5.3",                   "-0.72880726085203795",  "0.46710502065537601",   "173"),
    @(4, "2.7928199999999999",    "17.336165000000001",    "-9.6928970177535998",   "5.1177411031395398E-18","173"),
    @(5, "-0.109137188716422",    "0.10307401156550999",   "-1.4070165281359801",   "0.161216473982942",     "173"),
    @(6, "-2.0673605307369301E-2","1.9525071679182E-2",    "-0.33597637079760201",  "0.73729581626698404",   "173"),
    @(7, "19.027421176470501",    "17.419282222222201",    "1.3496483208698999",    "0.17889224836682099",   "173")
)

foreach ($r in $hypoRows) {
    $row = $r[0]
    $ws1.Cells.Item($row, 2).Value = $r[1]
    $ws1.Cells.Item($row, 3).Value = $r[2]
    $ws1.Cells.Item($row, 4).Value = $r[3]
    $ws1.Cells.Item($row, 5).Value = $r[4]
    $ws1.Cells.Item($row, 6).Value = $r[5]
}

# B:D -> 2 decimals, E -> 3 decimals, F (df) -> back to the plain default
$ws1.Range("B2:D7").NumberFormat = "0.00"
$ws1.Range("B9:C9").NumberFormat = "0.00"
$ws1.Range("E2:E7").NumberFormat = "0.000"
$ws1.Range("F2:F7").Style = "Normal"

$ws1.Range("B4").Select()

# ---------------------------------------------------------------------------
# 2. Add the "plex" worksheet after hypothesisTests
# ---------------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "plex"

$ws2.Range("B1").Value = "Control mean"
$ws2.Range("C1").Value = "Game mean"
$ws2.Range("D1").Value = "t"
$ws2.Range("E1").Value = "p"
$ws2.Range("F1").Value = "df"

$plexRows = @(
    @(2,  "Captivation",  "0.494117647058823",    "0.58888888888888802", "-1.2563065689892099", "0.21069866208711799", "173"),
    @(3,  "Challenge",    "0.34117647058823503",  "0.78888888888888797", "-6.6679140795516396",  "3.3519611511430401E-10","173"),
    @(4,  "Competition",  "0.152941176470588",    "0.64444444444444404", "-7.5994854731689996",  "1.7980839235228E-12", "173"),
    @(5,  "Completion",   "0.21176470588235199",  "0.62222222222222201", "-6.0040789378628503",  "1.1023184687861401E-8","173"),
    @(6,  "Discovery",    "0.78823529411764703",  "0.76666666666666605", "0.34081623040930897",  "0.73365546390720304", "173"),
    @(7,  "Progression",  "0.50588235294117601",  "0.74444444444444402", "-3.3492251800714499",  "9.9468405464890391E-4","173"),
    @(8,  "Exploration",  "0.90588235294117603",  "0.95555555555555505", "-1.2983661942134299",  "0.19588967143549901", "173"),
    @(9,  "Fantasy",      "0.51764705882352902",  "0.688888888888888",   "-2.3396479917131101",  "2.0443471097187998E-2","173"),
    @(10, "Humor",        "0.35294117647058798",  "0.27777777777777701", "1.0678612763083499",   "0.28707093300176301", "173"),
    @(11, "Nurture",      "0.48235294117646998",  "0.46666666666666601", "0.206535236271809",    "0.83661571731237006", "173"),
    @(12, "Relaxation",   "0.35294117647058798",  "0.51111111111111096", "-2.12494787487329",    "3.5010863573468901E-2","173"),
    @(13, "Sensation",    "0.51764705882352902",  "0.688888888888888",   "-2.3396479917131101",  "2.0443471097187998E-2","173")
)

foreach ($r in $plexRows) {
    $row = $r[0]
    $ws2.Cells.Item($row, 1).Value = $r[1]
    $ws2.Cells.Item($row, 2).Value = $r[2]
    $ws2.Cells.Item($row, 3).Value = $r[3]
    $ws2.Cells.Item($row, 4).Value = $r[4]
    $ws2.Cells.Item($row, 5).Value = $r[5]
    $ws2.Cells.Item($row, 6).Value = $r[6]
}

$ws2.Range("B1:F1").WrapText = $true
$ws2.Range("B2:D13").NumberFormat = "0.00"
$ws2.Range("E2:E13").NumberFormat = "0.000"
$ws2.Range("F2:F13").Style = "Normal"

$ws2.Range("D10").Select()

# Hand focus back to hypothesisTests (stays the active/selected tab)
$ws1.Activate()
$ws1.Range("B4").Select()
